$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric remain stored as text,
# matching the source data which uses inline text strings (e.g. European-style
# thousand separators, fixed 2-decimal prices, etc).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.682.57"
$ws.Range("E2").Value = "  +5.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.433.88"
$ws.Range("E3").Value = "  +7.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.42"
$ws.Range("E5").Value = "  +7.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.62"
$ws.Range("E6").Value = "  +7.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.440.87"
$ws.Range("E8").Value = "  +7.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.58"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("E11").Value = "  +8.51%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.023.13"
$ws.Range("E13").Value = "  +7.04%  "
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  +7.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.23"
$ws.Range("E16").Value = "  +5.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.772.36"
$ws.Range("E17").Value = "  +5.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.435.78"
$ws.Range("E18").Value = "  +7.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.27"
$ws.Range("E20").Value = "  +7.52%  "
$ws.Range("E21").Value = "  +3.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.19"
$ws.Range("E22").Value = "  +5.45%  "
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000107"
$ws.Range("E26").Value = "  +22.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  +10.19%  "
$ws.Range("E28").Value = "  +6.99%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.37"
$ws.Range("E32").Value = "  +15.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.76"
$ws.Range("E33").Value = "  +8.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.49"
$ws.Range("E34").Value = "  +4.49%  "
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.49"
$ws.Range("E37").Value = "  +8.93%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.05"
$ws.Range("E39").Value = "  +6.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0781"
$ws.Range("E40").Value = "  +10.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  +11.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.917.36"
$ws.Range("E42").Value = "  +4.36%  "
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.770"
$ws.Range("E44").Value = "  +7.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.85"
$ws.Range("E45").Value = "  +4.71%  "
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("E47").Value = "  +10.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.479.66"
$ws.Range("E48").Value = "  +7.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.60"
$ws.Range("E49").Value = "  +9.16%  "
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "295.73"
$ws.Range("E51").Value = "  +11.95%  "

# Row 23/24 swap (Polygon <-> Dai)
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.540"
$ws.Range("E24").Value = "  +3.10%  "

# Row 30/31 swap (PancakeSwap <-> RenderToken)
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.67"
$ws.Range("E30").Value = "  +8.60%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.04"
$ws.Range("E31").Value = "  +7.74%  "
